$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "64.603.53"
$ws.Range("E2").Value = "  +1.18%  "

# Row 3
$ws.Range("D3").Value = "3.098.43"
$ws.Range("E3").Value = "  -1.61%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.02"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +1.77%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.13"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.01%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.89"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.37%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.02"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.21%  "

# Row 8
$ws.Range("D8").Value = "3.093.15"
$ws.Range("E8").Value = "  -1.29%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.527"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.01%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.157"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.32%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.87"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.64%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.456"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.01%  "

# Row 13
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.43"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.67%  "

# Row 14
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000240"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.40%  "

# Row 15
$ws.Range("D15").Value = "3.710.88"
$ws.Range("E15").Value = "  +1.22%  "

# Row 16
$ws.Range("E16").Value = "  -1.95%  "

# Row 17
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "65.201.83"
$ws.Range("E17").Value = "  +2.49%  "

# Row 18
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.16"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.48%  "

# Row 19
$ws.Range("D19").Value = "3.124.81"
$ws.Range("E19").Value = "  -0.48%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "465.63"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.56%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.74"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.20%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.732"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.04%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.50"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.39%  "

# Row 24
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.12"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.77%  "

# Row 25
$ws.Range("B25").Value = "Fetch.AI"
$ws.Range("C25").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.35"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +6.27%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.33"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.50%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.38%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.69"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +4.64%  "

# Row 29
$ws.Range("B29").Value = "FirstDigitalUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.01"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.02%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.68"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.42%  "

# Row 31
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.25"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.24%  "

# Row 32
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.20"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.18%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.116"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +6.64%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.18"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.32%  "

# Row 35
$ws.Range("D35").Value = "0.0₃0845"
$ws.Range("E35").Value = "  -0.07%  "

# Row 36
$ws.Range("E36").Value = "  +1.67%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.09"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.73%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.33"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.13%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.25"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.32%  "

# Row 40
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.77"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.02%  "

# Row 41
$ws.Range("B41").Value = "Cosmos"
$ws.Range("C41").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.23"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +4.03%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "448.31"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.88%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.291"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +4.19%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0368"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.41%  "

# Row 45
$ws.Range("D45").Value = "2.852.37"
$ws.Range("E45").Value = "  -2.05%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.109"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.39%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.78"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.15%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.83"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.02%  "

# Row 49
$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.999"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.03%  "

# Row 50
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.93"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.80%  "

# Row 51
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.26"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.93%  "
